$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Evidence updates: four new Ryuk report rows (43-46) gain a
#     status/content tag and a report_date (column-at-a-time, matching
#     how the source data was merged in) ---
$ws.Range("A43").Value = "processed-8_34"
$ws.Range("A44").Value = "processed-17_30"
$ws.Range("G44").Value = "2020-11"
$ws.Range("G45").Value = "2020-10"
$ws.Range("A45").Value = "processed-17_26"
$ws.Range("A46").Value = "processed-19_25"
$ws.Range("G43").Value = "2021-01"
$ws.Range("G46").Value = "2020-10"

$ws.Range("B43").Value = "atk_ttp, htool, atk_chain, sigma, etpro"
$ws.Range("B44").Value = "atk_ttp, htool, atk_chain, sigma, etpro"
$ws.Range("B45").Value = "atk_ttp, htool, atk_chain, sigma, etpro"
$ws.Range("B46").Value = "atk_ttp, htool, atk_chain, sigma, etpro"

# --- Sheet refactor: widen/narrow the type & container columns ---
$ws.Columns.Item(3).ColumnWidth = 18.833333333333332
$ws.Columns.Item(4).ColumnWidth = 19.6640625

# --- Sheet refactor: move the active selection ---
$ws.Range("D26").Select()
